$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ between the old and new row contents.
$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "Q", "R", "S", "Z", "AB")

# Capture the current ("before") values of rows 2-5 for the columns above.
$before = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $before[$r] = $rowVals
}

# The edit is a cyclic rotation of the row data:
#   new row2 <- old row4
#   new row3 <- old row5
#   new row4 <- old row2
#   new row5 <- old row3
$mapping = @{ 2 = 4; 3 = 5; 4 = 2; 5 = 3 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
